{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that precedes it) that followed the\n// \"LOQ4213: Contabilidade e Custos (Requisito fraco)\" requirement line.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the three paragraphs to drop by their exact text content so the\n// edit is resilient to any surrounding content.\nconst targets = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the anchor paragraph (\"LOQ4213: ...\") and remove the block of\n// paragraphs that immediately follows it and matches `targets`.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"LOQ4213: Contabilidade e Custos (Requisito fraco)\") {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const start = anchorIndex + 1;\n  let matches = true;\n  for (let j = 0; j < targets.length; j++) {\n    const p = items[start + j];\n    if (!p || p.text !== targets[j]) {\n      matches = false;\n      break;\n    }\n  }\n  if (matches) {\n    for (let j = targets.length - 1; j >= 0; j--) {\n      items[start + j].delete();\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph that precedes it) that followed the\n# \"LOQ4213: Contabilidade e Custos (Requisito fraco)\" requirement line.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"LOQ4213: Contabilidade e Custos (Requisito fraco)\"\n$targetTexts = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    [char]0xA9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ne -1) {\n    $startIndex = $anchorIndex + 1\n    $matches = $true\n    for ($j = 0; $j -lt $targetTexts.Length; $j++) {\n        $idx = $startIndex + $j\n        if ($idx -gt $d.Paragraphs.Count) {\n            $matches = $false\n            break\n        }\n        $t = $d.Paragraphs.Item($idx).Range.Text.TrimEnd(\"`r\", \"`a\")\n        if ($t -ne $targetTexts[$j]) {\n            $matches = $false\n            break\n        }\n    }\n\n    if ($matches) {\n        $firstPara = $d.Paragraphs.Item($startIndex)\n        $lastPara = $d.Paragraphs.Item($startIndex + $targetTexts.Length - 1)\n        $r = $d.Range($firstPara.Range.Start, $lastPara.Range.End)\n        $r.Delete()\n    }\n}\n"}
